$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collect all D-column cells that need new text values so we can force
# text storage (avoiding Excel's automatic number conversion) and then
# strip the temporary text number-format back off so styling is untouched.
$dRows = @(2,3,5,6,8,9,12,13,14,15,16,17,18,20,21,22,23,24,25,26,28,29,30,31,32,33,36,37,38,40,41,42,43,44,45,46,49,50,51)
foreach ($r in $dRows) { $ws.Cells.Item($r, 4).NumberFormat = "@" }

$ws.Range("D2").Value = "29.229.88"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.861.35"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "0.7115"
$ws.Range("E5").Value = "  +1.19%  "
$ws.Range("D6").Value = "237.94"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.08182"
$ws.Range("E8").Value = "  +10.47%  "
$ws.Range("D9").Value = "0.3043"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  -0.88%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("D12").Value = "1.885.13"
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "5.175"
$ws.Range("E13").Value = "  -0.80%  "
$ws.Range("D14").Value = "0.7094"
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("D15").Value = "89.67"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "29.270.70"
$ws.Range("E16").Value = "  +0.35%  "
$ws.Range("D17").Value = "0.000007927"
$ws.Range("E17").Value = "  +3.67%  "
$ws.Range("D18").Value = "5.794"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").Value = "  +2.20%  "
$ws.Range("D20").Value = "237.35"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "2.108.70"
$ws.Range("E22").Value = "  +0.70%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "7.428"
$ws.Range("E24").Value = "  -2.27%  "
$ws.Range("D25").Value = "162.75"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("D26").Value = "8.968"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "18.10"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Value = "1.961"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("D30").Value = "1.425"
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("D31").Value = "1.484"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("D32").Value = "4.398"
$ws.Range("E32").Value = "  -2.63%  "
$ws.Range("D33").Value = "4.027"
$ws.Range("E33").Value = "  +0.92%  "
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").Value = "0.7092"
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("D38").Value = "2.672"
$ws.Range("E38").Value = "  +0.22%  "
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("D40").Value = "2.732"
$ws.Range("E40").Value = "  +2.24%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.9245"
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.138.57"
$ws.Range("E42").Value = "  +6.21%  "
$ws.Range("D43").Value = "0.4284"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "5.903"
$ws.Range("E44").Value = "  -1.68%  "
$ws.Range("D45").Value = "70.41"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("D49").Value = "2.010.00"
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").Value = "9.222"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").Value = "6.951"
$ws.Range("E51").Value = "  -1.20%  "

# Remove the temporary number-format override, restoring the cells to
# their original (unstyled) appearance while keeping the text values.
foreach ($r in $dRows) { $ws.Cells.Item($r, 4).ClearFormats() }
